$d = $word.ActiveDocument

# The paragraph under "MetaModel / Backend:" / "Reference Model / Occurrences
# annotations matrix:" reads:
#   "Attachment: <hyperlink to the scanned PDF> shows an example of a fully
#    expanded set ..."
# The edit drops the hyperlink (and its URL text) entirely and folds the
# remaining sentence so it now reads "Attachment: show an example ...".

if ($d.Hyperlinks.Count -gt 0) {
    $h = $d.Hyperlinks.Item(1)
    $hlRange = $d.Range($h.Range.Start, $h.Range.End)
    $hlRange.Text = ""
}

$d.Content.Find.Execute(" shows an example", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "show an example", 2)
